$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A,B,C=15.83.., D,E=20.83.., F=15.83.., G=25.83..) ---
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 20
$ws.Columns.Item(5).ColumnWidth = 20
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(7).ColumnWidth = 25

# --- Make sure every cell in the new used range keeps a plain "General"
#     style (matches original s="1"), then force column A to Text so the
#     dd.mm.yyyy-looking strings are not reinterpreted as date serials. ---
$ws.Range("A1:G5").NumberFormat = "General"
$ws.Range("A1:A5").NumberFormat = "@"

# --- Row 1: header (now capitalised) ---
$ws.Range("A1").Value = "Дата"
$ws.Range("B1").Value = "Сумма"
$ws.Range("C1").Value = "Вид"
$ws.Range("D1").Value = "Люди"
$ws.Range("E1").Value = "Титул"
$ws.Range("F1").Value = "Объект"
$ws.Range("G1").Value = "Бригады"

# --- Row 2 (new content) ---
$ws.Range("A2").Value = "10.11.2022"
$ws.Range("B2").Value = "253,6"
$ws.Range("C2").Value = "Бригада"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Аванс"
$ws.Range("F2").Value = "MCM project"
$ws.Range("G2").Value = "Shartukh Anton i Aliaksandr"

# --- Row 3 (new content) ---
$ws.Range("A3").Value = "10.11.2022"
$ws.Range("B3").Value = "23,69"
$ws.Range("C3").Value = "Люди"
$ws.Range("D3").Value = "Tarasiuk Oleksandr"
$ws.Range("E3").Value = "Страховки на авто"
$ws.Range("F3").Value = "Skysawa"
$ws.Range("G3").Value = ""

# --- Row 4 (old row 2 data, unchanged) ---
$ws.Range("A4").Value = "10.11.2022"
$ws.Range("B4").Value = "25,36"
$ws.Range("C4").Value = "Общее"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "Зарплата"
$ws.Range("F4").Value = "Karpacz"
$ws.Range("G4").Value = ""

# --- Row 5 (old row 3 data, with updates) ---
$ws.Range("A5").Value = "10.11.2022"
$ws.Range("B5").Value = "25,69"
$ws.Range("C5").Value = "Люди"
$ws.Range("D5").Value = "Анастасия PM"
$ws.Range("E5").Value = "Топливо"
$ws.Range("F5").Value = "Office"
$ws.Range("G5").Value = ""
